$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 53, shifting existing rows 53-135 down to 54-136.
$ws.Rows(53).Insert()

# Populate the newly inserted row 53 with the new data record.
$ws.Range("A53").Value = 8
$ws.Range("B53").Value = "Terminal La Palmera de La Serena"
$ws.Range("C53").Value = "Coquimbo"
$ws.Range("D53").Value = 45079
$ws.Range("E53").Value = 4
$ws.Range("F53").Value = 100114007
$ws.Range("G53").Value = "Jengibre"
$ws.Range("H53").Value = "Sin especificar"
$ws.Range("I53").Value = "Primera"
$ws.Range("J53").Value = 400
$ws.Range("K53").Value = 17000
$ws.Range("L53").Value = 18000
$ws.Range("M53").Value = 17500
$ws.Range("N53").Value = "$/caja 13 kilos"
$ws.Range("O53").Value = "Perú"
$ws.Range("P53").Value = 1346
$ws.Range("Q53").Value = 13
$ws.Range("R53").Value = "Hortaliza"
